$d = $word.ActiveDocument
$p = $d.Paragraphs(1)
$r = $p.Range
$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pBdr><w:top w:space='5'/><w:left w:space='5'/><w:bottom w:space='5'/><w:right w:space='5'/></w:pBdr><w:spacing w:after='0'/><w:ind w:left='225'/><w:jc w:val='left'/></w:pPr><w:r><w:rPr><w:rFonts w:ascii='Times New Roman' w:hAnsi='Times New Roman'/><w:b w:val='false'/><w:i w:val='false'/><w:color w:val='000000'/><w:sz w:val='22'/></w:rPr><w:t>**ID__AFFARS_SUBPART_5342_2__ID**</w:t></w:r></w:p>"
$r.InsertXML($xml)
Write-Host "Para count after:" $d.Paragraphs.Count
